# Region 4B.xlsx edit:
#  - Fix a typo in the second sheet's name: "Puerto Prinsesa" -> "Puerto Princesa"
#  - Restore the cursor/selection on that sheet to cell D17 (matches the
#    author's last saved selection before committing the file).

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Puerto Prinsesa")
$ws2.Name = "Puerto Princesa"

$ws2.Activate()
$ws2.Range("D17").Select() | Out-Null
